# teste 01 - modelagem parametrica
#
# The "BIM 3D - Modelagem Paramétrica" course (row 7) got its end date
# updated and a grade entered for Unidade 1. The following row, "BIM 4D -
# Planejamento e Controle de Obras" (row 8), also had its end date updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Row 8: BIM 4D - Planejamento e Controle de Obras
$ws.Range("D8").Value = "29/10/2020"

# Row 7: BIM 3D - Modelagem Paramétrica
$ws.Range("D7").Value = "29/08/2020"
$ws.Range("F7").Value = 10

# Reflect where the user ended up working in the sheet.
$ws.Range("D4").Select()
